$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix letter case inconsistencies in header labels (SSDM-12286)
# "Vocabulary Code" -> "Vocabulary code" (property type header, row 4)
$ws.Range("H4").Value = "Vocabulary code"
# "Generated Code Prefix" -> "Generated code prefix" (sample type header, row 2)
$ws.Range("E2").Value = "Generated code prefix"

# Leave the active selection on the last edited cell
$ws.Range("E2").Select()
